# Update the as_of_utc timestamps (column AA, rows 2-26) on the
# "Главные" and "Линейные" sheets from 2025-11-09 03:02:37 to 2025-11-09 04:53:53.

$wb = $excel.ActiveWorkbook

$sheetNames = @("Главные", "Линейные")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    for ($row = 2; $row -le 26; $row++) {
        $ws.Range("AA$row").Value = "2025-11-09 04:53:53"
    }
}
